$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '57.896.34'
Set-TextValue 'E2' '  -0.38%  '

# Row 3
Set-TextValue 'D3' '2.349.96'
Set-TextValue 'E3' '  +0.74%  '

# Row 4
Set-TextValue 'E4' '  +0.06%  '

# Row 5
Set-TextValue 'D5' '540.65'
Set-TextValue 'E5' '  -0.38%  '

# Row 6
Set-TextValue 'D6' '134.62'
Set-TextValue 'E6' '  -0.20%  '

# Row 7
Set-TextValue 'E7' '  +0.39%  '

# Row 8
Set-TextValue 'D8' '0.571'
Set-TextValue 'E8' '  +6.45%  '

# Row 9
Set-TextValue 'D9' '0.103'
Set-TextValue 'E9' '  +0.52%  '

# Row 10
Set-TextValue 'D10' '5.55'
Set-TextValue 'E10' '  +2.34%  '

# Row 12
Set-TextValue 'D12' '0.356'
Set-TextValue 'E12' '  +0.79%  '

# Row 13
Set-TextValue 'D13' '2.769.30'
Set-TextValue 'E13' '  +0.30%  '

# Row 14
Set-TextValue 'D14' '23.74'
Set-TextValue 'E14' '  +0.89%  '

# Row 15
Set-TextValue 'D15' '57.831.03'
Set-TextValue 'E15' '  -0.27%  '

# Row 16
Set-TextValue 'D16' '0.0000135'
Set-TextValue 'E16' '  +0.61%  '

# Row 17
Set-TextValue 'D17' '2.321.46'
Set-TextValue 'E17' '  -1.44%  '

# Row 18
Set-TextValue 'D18' '10.70'
Set-TextValue 'E18' '  +1.10%  '

# Row 19
Set-TextValue 'B19' 'Polkadot'
Set-TextValue 'C19' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D19' '4.29'
Set-TextValue 'E19' '  +1.28%  '

# Row 20
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '329.99'
Set-TextValue 'E20' '  -2.48%  '

# Row 21
Set-TextValue 'D21' '6.72'
Set-TextValue 'E21' '  -2.15%  '

# Row 22
Set-TextValue 'E22' '  +0.16%  '

# Row 23
Set-TextValue 'D23' '62.66'
Set-TextValue 'E23' '  +0.89%  '

# Row 24
Set-TextValue 'E24' '  -1.96%  '

# Row 25
Set-TextValue 'D25' '0.998'
Set-TextValue 'E25' '  -0.16%  '

# Row 26
Set-TextValue 'D26' '8.36'
Set-TextValue 'E26' '  -1.61%  '

# Row 27
Set-TextValue 'D27' '1.35'
Set-TextValue 'E27' '  -5.52%  '

# Row 28
Set-TextValue 'E28' '  +0.08%  '

# Row 29
Set-TextValue 'D29' '170.12'
Set-TextValue 'E29' '  -0.82%  '

# Row 30
Set-TextValue 'E30' '  -0.35%  '

# Row 31
Set-TextValue 'E31' '  -1.24%  '

# Row 32
Set-TextValue 'E32' '  +0.46%  '

# Row 33
Set-TextValue 'E33' '  -1.19%  '

# Row 34
Set-TextValue 'E34' '  +0.00%  '

# Row 35
Set-TextValue 'E35' '  +0.33%  '

# Row 36
Set-TextValue 'E36' '  +1.04%  '

# Row 37
Set-TextValue 'E37' '  -2.02%  '

# Row 38
Set-TextValue 'E38' '  -0.33%  '

# Row 39
Set-TextValue 'D39' '39.07'
Set-TextValue 'E39' '  -0.77%  '

# Row 40
Set-TextValue 'D40' '142.81'
Set-TextValue 'E40' '  -4.30%  '

# Row 41
Set-TextValue 'D41' '0.379'
Set-TextValue 'E41' '  +0.19%  '

# Row 42
Set-TextValue 'E42' '  +0.44%  '

# Row 43
Set-TextValue 'D43' '288.43'
Set-TextValue 'E43' '  +0.82%  '

# Row 44
Set-TextValue 'E44' '  +1.77%  '

# Row 45
Set-TextValue 'E45' '  +0.51%  '

# Row 46
Set-TextValue 'D46' '19.15'
Set-TextValue 'E46' '  -0.69%  '

# Row 47
Set-TextValue 'D47' '0.566'
Set-TextValue 'E47' '  +0.93%  '

# Row 48
Set-TextValue 'E48' '  +1.67%  '

# Row 49
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '17.43'
Set-TextValue 'E49' '  -0.83%  '

# Row 50
Set-TextValue 'B50' 'Polygon'
Set-TextValue 'C50' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D50' '0.379'
Set-TextValue 'E50' '  -1.01%  '

# Row 51
Set-TextValue 'E51' '  +0.58%  '
